$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextValue "D2" "60.645.16"
Set-TextValue "E2" "  +0.83%  "
Set-TextValue "D3" "2.599.67"
Set-TextValue "E3" "  +0.58%  "
Set-TextValue "E4" "  -0.13%  "
Set-TextValue "D5" "517.62"
Set-TextValue "E5" "  +2.34%  "
Set-TextValue "D6" "154.15"
Set-TextValue "E6" "  +0.59%  "
Set-TextValue "E7" "  -0.16%  "
Set-TextValue "E8" "  +3.27%  "
Set-TextValue "D9" "6.71"
Set-TextValue "E9" "  +1.14%  "
Set-TextValue "D10" "0.106"
Set-TextValue "E10" "  +2.58%  "
Set-TextValue "D11" "0.347"
Set-TextValue "E11" "  +0.09%  "
Set-TextValue "D12" "0.130"
Set-TextValue "E12" "  +1.47%  "
Set-TextValue "D13" "3.050.83"
Set-TextValue "E13" "  +0.24%  "
Set-TextValue "D14" "60.672.35"
Set-TextValue "E14" "  +0.70%  "
Set-TextValue "D15" "21.73"
Set-TextValue "E15" "  +0.90%  "
Set-TextValue "D16" "0.0000141"
Set-TextValue "E16" "  +0.77%  "
Set-TextValue "D17" "2.602.19"
Set-TextValue "E17" "  +0.21%  "
Set-TextValue "D18" "4.75"
Set-TextValue "E18" "  -1.62%  "
Set-TextValue "D19" "352.24"
Set-TextValue "E19" "  +1.55%  "
Set-TextValue "D20" "10.57"
Set-TextValue "E20" "  +2.10%  "
Set-TextValue "D21" "6.22"
Set-TextValue "E21" "  +1.50%  "
Set-TextValue "D22" "1.00"
Set-TextValue "E22" "  +0.10%  "
Set-TextValue "D23" "61.03"
Set-TextValue "E23" "  +1.59%  "
Set-TextValue "D24" "0.427"
Set-TextValue "E24" "  +1.43%  "
Set-TextValue "E25" "  -0.02%  "
Set-TextValue "D26" "2.711.59"
Set-TextValue "E26" "  +0.30%  "
Set-TextValue "D27" "0.999"
Set-TextValue "E27" "  +0.02%  "
Set-TextValue "D28" "0.0₃0843"
Set-TextValue "E28" "  -0.07%  "
Set-TextValue "D29" "7.34"
Set-TextValue "E29" "  -0.94%  "
Set-TextValue "E30" "  -0.11%  "
Set-TextValue "D31" "6.28"
Set-TextValue "E31" "  +9.35%  "
Set-TextValue "D32" "19.39"
Set-TextValue "E32" "  +0.48%  "
Set-TextValue "D33" "1.60"
Set-TextValue "E33" "  +2.68%  "
Set-TextValue "D34" "149.74"
Set-TextValue "E34" "  -3.00%  "
Set-TextValue "D35" "4.19"
Set-TextValue "E35" "  +5.18%  "
Set-TextValue "D36" "0.925"
Set-TextValue "E36" "  +10.32%  "
Set-TextValue "D37" "1.20"
Set-TextValue "E37" "  +0.94%  "
Set-TextValue "D38" "1.49"
Set-TextValue "E38" "  +1.93%  "
Set-TextValue "D39" "3.79"
Set-TextValue "E39" "  +0.75%  "
Set-TextValue "D40" "36.37"
Set-TextValue "E40" "  +1.59%  "
Set-TextValue "D41" "0.843"
Set-TextValue "E41" "  -0.37%  "
Set-TextValue "D42" "287.18"
Set-TextValue "E42" "  -3.25%  "
Set-TextValue "E43" "  +1.42%  "
Set-TextValue "D44" "0.623"
Set-TextValue "E44" "  +1.23%  "
Set-TextValue "D45" "0.0560"
Set-TextValue "E45" "  +0.13%  "
Set-TextValue "E46" "  -0.08%  "
Set-TextValue "D47" "19.55"
Set-TextValue "E47" "  -0.93%  "
Set-TextValue "B48" "RenderToken"
Set-TextValue "C48" "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue "D48" "4.95"
Set-TextValue "E48" "  +1.24%  "
Set-TextValue "B49" "VeChain"
Set-TextValue "C49" "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue "D49" "0.0237"
Set-TextValue "E49" "  +1.36%  "
Set-TextValue "E50" "  +0.12%  "
Set-TextValue "D51" "19.14"
Set-TextValue "E51" "  +8.54%  "
